$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing match rows (old rows 5 and 6); row 4 is refreshed
# in-place below with an updated odds snapshot for the same fixture.
$ws.Rows("5:6").Delete()

# Refresh row 4 with the latest odds snapshot (Coritiba vs Botafogo SP).
$ws.Range("A4").Value = "bcwMI4Uc"
$ws.Range("B4").Value = "22/11/2024"
$ws.Range("C4").Value = "19:00"
$ws.Range("D4").Value = "BRAZIL - SERIE B"
$ws.Range("E4").Value = "Coritiba"
$ws.Range("F4").Value = "Botafogo SP"
$ws.Range("G4").Value = 1.55
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 6.5
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 2.2
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 5.5
$ws.Range("X4").Value = 6.5
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 7.5
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 81
$ws.Range("AG4").Value = ""
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 67
$ws.Range("AL4").Value = 51
$ws.Range("AM4").Value = 51
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 8
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 26
$ws.Range("AR4").Value = 51
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.63
$ws.Range("AU4").Value = 9.5
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 7.5
$ws.Range("AX4").Value = 34
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 151
$ws.Range("BA4").Value = 201
$ws.Range("BB4").Value = 81
$ws.Range("BC4").Value = 81
$ws.Range("BD4").Value = 81
